{"js": "// Update the date in the title paragraph.\nconst titleResults = context.document.body.search(\"2024-09-27 Friday\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"2024-09-28 Saturday\", \"Replace\");\n}\nawait context.sync();\n\n// Update the practice-problem answers inside the single table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Only the 1st, 5th, 9th, 13th and 17th rows (0-based: 0, 4, 8, 12, 16)\n// contain text; the rows in between are intentionally blank spacer rows.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"99\u00f77=14, 1\", newText: \"20\u00f72=10, 0\" },\n  { row: 0, col: 1, oldText: \"17\u00f77=2, 3\", newText: \"44\u00f73=14, 2\" },\n  { row: 0, col: 2, oldText: \"67\u00f79=7, 4\", newText: \"10\u00f75=2, 0\" },\n  { row: 0, col: 3, oldText: \"71\u00f78=8, 7\", newText: \"20\u00f72=10, 0\" },\n  { row: 0, col: 4, oldText: \"32\u00f73=10, 2\", newText: \"71\u00f76=11, 5\" },\n\n  { row: 4, col: 0, oldText: \"10\u00f75=2, 0\", newText: \"26\u00f77=3, 5\" },\n  { row: 4, col: 1, oldText: \"69\u00f79=7, 6\", newText: \"77\u00f77=11, 0\" },\n  { row: 4, col: 2, oldText: \"54\u00f75=10, 4\", newText: \"82\u00f73=27, 1\" },\n  { row: 4, col: 3, oldText: \"99\u00f73=33, 0\", newText: \"10\u00f76=1, 4\" },\n  { row: 4, col: 4, oldText: \"95\u00f79=10, 5\", newText: \"79\u00f77=11, 2\" },\n\n  { row: 8, col: 0, oldText: \"26\u00f74=6, 2\", newText: \"71\u00f79=7, 8\" },\n  { row: 8, col: 1, oldText: \"58\u00f78=7, 2\", newText: \"56\u00f72=28, 0\" },\n  { row: 8, col: 2, oldText: \"63\u00f75=12, 3\", newText: \"91\u00f79=10, 1\" },\n  { row: 8, col: 3, oldText: \"81\u00f75=16, 1\", newText: \"53\u00f76=8, 5\" },\n  { row: 8, col: 4, oldText: \"57\u00f79=6, 3\", newText: \"34\u00f75=6, 4\" },\n\n  { row: 12, col: 0, oldText: \"87\u00f75=17, 2\", newText: \"64\u00f72=32, 0\" },\n  { row: 12, col: 1, oldText: \"64\u00f74=16, 0\", newText: \"11\u00f75=2, 1\" },\n  { row: 12, col: 2, oldText: \"19\u00f79=2, 1\", newText: \"57\u00f74=14, 1\" },\n  { row: 12, col: 3, oldText: \"25\u00f72=12, 1\", newText: \"60\u00f79=6, 6\" },\n  { row: 12, col: 4, oldText: \"19\u00f79=2, 1\", newText: \"86\u00f78=10, 6\" },\n\n  { row: 16, col: 0, oldText: \"76\u00f74=19, 0\", newText: \"22\u00f72=11, 0\" },\n  { row: 16, col: 1, oldText: \"94\u00f77=13, 3\", newText: \"38\u00f76=6, 2\" },\n  { row: 16, col: 2, oldText: \"22\u00f74=5, 2\", newText: \"61\u00f77=8, 5\" },\n  { row: 16, col: 3, oldText: \"18\u00f75=3, 3\", newText: \"88\u00f79=9, 7\" },\n  { row: 16, col: 4, oldText: \"13\u00f72=6, 1\", newText: \"13\u00f73=4, 1\" },\n];\n\nconst allResults = [];\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const found = cell.body.search(r.oldText, { matchCase: true });\n  found.load(\"items\");\n  allResults.push(found);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const found = allResults[i];\n  if (found.items.length > 0) {\n    found.items[0].insertText(replacements[i].newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date in the title paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2024-09-28 Saturday\"\n\n# Update the practice-problem answers inside the single table.\n$tbl = $d.Tables.Item(1)\n\n# Word COM table rows/columns are 1-based. Only rows 1, 5, 9, 13 and 17\n# contain text; the rows in between are intentionally blank spacer rows.\n$tbl.Cell(1, 1).Range.Text = \"20\u00f72=10, 0\"\n$tbl.Cell(1, 2).Range.Text = \"44\u00f73=14, 2\"\n$tbl.Cell(1, 3).Range.Text = \"10\u00f75=2, 0\"\n$tbl.Cell(1, 4).Range.Text = \"20\u00f72=10, 0\"\n$tbl.Cell(1, 5).Range.Text = \"71\u00f76=11, 5\"\n\n$tbl.Cell(5, 1).Range.Text = \"26\u00f77=3, 5\"\n$tbl.Cell(5, 2).Range.Text = \"77\u00f77=11, 0\"\n$tbl.Cell(5, 3).Range.Text = \"82\u00f73=27, 1\"\n$tbl.Cell(5, 4).Range.Text = \"10\u00f76=1, 4\"\n$tbl.Cell(5, 5).Range.Text = \"79\u00f77=11, 2\"\n\n$tbl.Cell(9, 1).Range.Text = \"71\u00f79=7, 8\"\n$tbl.Cell(9, 2).Range.Text = \"56\u00f72=28, 0\"\n$tbl.Cell(9, 3).Range.Text = \"91\u00f79=10, 1\"\n$tbl.Cell(9, 4).Range.Text = \"53\u00f76=8, 5\"\n$tbl.Cell(9, 5).Range.Text = \"34\u00f75=6, 4\"\n\n$tbl.Cell(13, 1).Range.Text = \"64\u00f72=32, 0\"\n$tbl.Cell(13, 2).Range.Text = \"11\u00f75=2, 1\"\n$tbl.Cell(13, 3).Range.Text = \"57\u00f74=14, 1\"\n$tbl.Cell(13, 4).Range.Text = \"60\u00f79=6, 6\"\n$tbl.Cell(13, 5).Range.Text = \"86\u00f78=10, 6\"\n\n$tbl.Cell(17, 1).Range.Text = \"22\u00f72=11, 0\"\n$tbl.Cell(17, 2).Range.Text = \"38\u00f76=6, 2\"\n$tbl.Cell(17, 3).Range.Text = \"61\u00f77=8, 5\"\n$tbl.Cell(17, 4).Range.Text = \"88\u00f79=9, 7\"\n$tbl.Cell(17, 5).Range.Text = \"13\u00f73=4, 1\"\n"}
